# Add thêm nhân sự Nguyễn Hữu Quang
# Updates the "Lương" (Salary) sheet: staff count, total work days,
# base salary and the resulting totals that roll up from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B1").Value = 3
$ws.Range("B2").Value = 29
$ws.Range("B3").Value = 3107142.857142857
$ws.Range("B31").Value = 3307142.857142857
$ws.Range("B34").Value = 3307142.857142857
